# Clean up the price tracker sheet ("add readme and clean data"):
#  - drop the ASIN / Bar Code / Case Qty columns (C, D, E)
#  - keep only Product Description / Cost Price (columns A, B)
#  - replace the old product rows with the current product list
#  - add two new products: Apple Airpods Pro 3, Samsung Galaxy S25

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the now-unused ASIN / Bar Code / Case Qty columns entirely.
$ws.Range("C1:E4").Clear()

# Header row.
$ws.Range("A1").Value = "Product Description"
$ws.Range("B1").Value = " Cost Price"
$ws.Range("B1").Style = "Normal"

# Product rows.
$ws.Range("A2").Value = "Iphone 17 pro"
$ws.Range("B2").Value = 134000
$ws.Range("B2").Style = "Normal"

$ws.Range("A3").Value = "Google Pixel 9"
$ws.Range("B3").Value = 64000
$ws.Range("B3").Style = "Normal"

$ws.Range("A4").Value = "Google Pixel 9a"
$ws.Range("B4").Value = 43000
$ws.Range("B4").Style = "Normal"

$ws.Range("A5").Value = "Apple Airpods Pro 3"
$ws.Range("B5").Value = 25000

$ws.Range("A6").Value = "Samsung Galaxy S25"
$ws.Range("B6").Value = 75000

# Column A needs to widen to fit the longer product names.
$ws.Columns.Item(1).ColumnWidth = 17.8

$ws.Range("C8").Select()
